# ---------------------------------------------------------------------------
# Updates workbook "STROBEL CORDERO MARIA ELISABETH.xlsx":
#   1. Sheet "VENTAS POR GRUPO": insert a new client row ("GAVILANES VELEZ
#      MARIA VALERIA", under advisor "OFICINA-CATAECSA") before the existing
#      row for "GRANJA VANEGAS MARCELA" (currently row 292), pushing all
#      following rows down by one. Also bump the "N de 341" summary labels
#      on the (now shifted) totals row to "N de 342", and fix the count for
#      column H (8 -> 7) to reflect the H221 edit below. Two independent
#      value corrections: M82 and H221.
#   2. Sheet "VENTA MENSUAL": same new-client row insertion (currently before
#      row 296), plus independent value corrections F82 and F225, plus the
#      resulting change to the "octubre" (F) grand total.
#   3. Sheet "CUMPLIMIENTO MENSUAL": update the precomputed VENTA / POR
#      CUMPLIR / CUMPLIMIENTO figures on the rows affected by the above
#      corrections (rows 24, 54, 63) and the TOTAL row (76).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1) VENTAS POR GRUPO
# ===========================================================================
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# --- insert the new client row at 292 (shifts 292..343 down to 293..344) ---
$ws1.Rows.Item(292).Insert()

$ws1.Cells.Item(292, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(292, 2).Value = "GAVILANES VELEZ MARIA VALERIA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(292, $c).Value = 0
}

# --- two standalone value corrections (rows before the insertion point,
#     so their row numbers are unaffected by the shift) ---
$ws1.Cells.Item(82, 13).Value = 250.09    # M82: PORCELANATO, was 111.13
$ws1.Cells.Item(221, 8).Value = 0         # H221: INODOROS, was 71.09999999999999

# --- refresh the "N de 341" -> "N de 342" labels on the totals row, now
#     shifted from row 343 to row 344; column H's count also drops 8 -> 7
#     because the H221 correction above zeroed out a previously non-zero
#     value ---
$ws1.Cells.Item(344, 3).Value  = "5 de 342"
$ws1.Cells.Item(344, 4).Value  = "15 de 342"
$ws1.Cells.Item(344, 5).Value  = "7 de 342"
$ws1.Cells.Item(344, 6).Value  = "0 de 342"
$ws1.Cells.Item(344, 7).Value  = "0 de 342"
$ws1.Cells.Item(344, 8).Value  = "7 de 342"
$ws1.Cells.Item(344, 9).Value  = "10 de 342"
$ws1.Cells.Item(344, 10).Value = "1 de 342"
$ws1.Cells.Item(344, 11).Value = "1 de 342"
$ws1.Cells.Item(344, 12).Value = "5 de 342"
$ws1.Cells.Item(344, 13).Value = "29 de 342"
$ws1.Cells.Item(344, 14).Value = "0 de 342"
$ws1.Cells.Item(344, 15).Value = "0 de 342"
$ws1.Cells.Item(344, 16).Value = "3 de 342"
$ws1.Cells.Item(344, 17).Value = "0 de 342"
$ws1.Cells.Item(344, 18).Value = "0 de 342"

# ===========================================================================
# 2) VENTA MENSUAL
# ===========================================================================
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# --- insert the matching new client row at 296 ---
$ws2.Rows.Item(296).Insert()

$ws2.Cells.Item(296, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(296, 2).Value = "GAVILANES VELEZ MARIA VALERIA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(296, $c).Value = 0
}

# --- two standalone value corrections (unaffected by the shift, both rows
#     are above the insertion point) ---
$ws2.Cells.Item(82, 6).Value  = 640.46   # F82: octubre, was 501.5
$ws2.Cells.Item(225, 6).Value = 0        # F225: octubre, was 71.09999999999999

# --- grand-total row, now shifted from 347 to 348: only the "octubre"
#     column moves, by the net of the two corrections above ---
$ws2.Cells.Item(348, 6).Value = 88622.72   # was 88554.86

# ===========================================================================
# 3) CUMPLIMIENTO MENSUAL
# ===========================================================================
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 24: CASTRO ALCIVAR EDA MARIA / PORCELANATO
$ws3.Cells.Item(24, 4).Value = 5227.26               # VENTA
$ws3.Cells.Item(24, 5).Value = 43396.8                # POR CUMPLIR
$ws3.Cells.Item(24, 6).Value = 0.1075035692206698     # CUMPLIMIENTO

# Row 54: LINDAO ZUÑIGA BRYAN JOSE / INODOROS
$ws3.Cells.Item(54, 4).Value = 394.37
$ws3.Cells.Item(54, 5).Value = 455.47419682004
$ws3.Cells.Item(54, 6).Value = 0.4640497652106818

# Row 63: OFICINA-CATAECSA / PORCELANATO
$ws3.Cells.Item(63, 4).Value = 5264.61
$ws3.Cells.Item(63, 5).Value = 14735.39
$ws3.Cells.Item(63, 6).Value = 0.2632305

# Row 76: TOTAL
$ws3.Cells.Item(76, 4).Value = 87742.86
$ws3.Cells.Item(76, 5).Value = 319869.0970193434
$ws3.Cells.Item(76, 6).Value = 0.2152607608511252
